$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
#
#    We copy the existing bold "Play Black Hole Slot..." paragraph from
#    near the end of the document (it already has the same run layout
#    we need: an empty leading run followed by a bold run) and paste it
#    in place, then edit its text. This preserves the empty leading
#    <w:r/> run that the rest of the document's paragraphs use.
# ------------------------------------------------------------------

$countBefore = $d.Paragraphs.Count
$boldTemplatePara = $d.Paragraphs.Item($countBefore - 1)
$boldTemplatePara.Range.Copy()

$titlePara = $d.Paragraphs.Item(1)
$titleEnd = $titlePara.Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

$metaStart = $metaPara.Range.Start
$metaEnd = $metaPara.Range.End

# Replace the pasted (bold) run's text with "Meta description"
$boldRng = $d.Range($metaStart, $metaEnd - 1)
$boldRng.Text = "Meta description"

# Append the remaining (non-bold) text of the meta description
$metaParaNow = $d.Paragraphs.Item(2)
$metaParaNow.Range.InsertAfter(": Read our review of Black Hole slot game with its unique gameplay structure, impressive graphics, and distinctive black hole feature. Play for free now.")

# ------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Black Hole Slot..." paragraph that
#    used to sit right before the final (italic) paragraph.
# ------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupBoldPara = $d.Paragraphs.Item($count - 1)
$dupBoldPara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the DALLE
#    image-generation prompt, keeping the existing run/format layout.
# ------------------------------------------------------------------

$countAfter = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($countAfter)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End

$lastTextRng = $d.Range($lastStart, $lastEnd - 1)
$lastTextRng.Text = 'Prompt for DALLE: Create a cartoon-style feature image for the game "Black Hole" that features a happy Maya warrior wearing glasses. The image should have a background of space with black holes, paying homage to the game''s theme, with the Maya warrior standing in front of the game reels with fruits, stars, bells, and number sevens visible on the screen. The image should be colorful and eye-catching, with the Maya warrior looking excited and ready to play the game.'
